# CDM_Floor.xlsx edit: reorder/rename header columns, drop the now-unused
# measurement columns (I:M), remove the stray formatted blank row 20, resize
# the remaining data columns to fit their new headers, and move the active
# selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers: new column order / names -------------------------------
# A:FloorId  B:FloorCode  C:Name  D:FloorNumber  E:ValidFrom  F:ValidUntil
# G:BuildingId  H:Guid (renamed from GUID)
$ws.Range("A1").Value = "FloorId"
$ws.Range("B1").Value = "FloorCode"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "FloorNumber"
$ws.Range("E1").Value = "ValidFrom"
$ws.Range("F1").Value = "ValidUntil"
$ws.Range("G1").Value = "BuildingId"
$ws.Range("H1").Value = "Guid"

# Drop the old AreaMeasurementId/FloorHightRaw/FloorLevelShell/
# ConstructionAreaBearing/ConstructionAreaNonBearing columns (I:M) — no
# longer part of the sheet.
$ws.Range("I1:M1").Clear()

# Remove the stray formatted-but-empty row 20 (A20 only carried a style).
$ws.Range("A20").Clear()

# --- Column widths: fit the surviving data columns to their headers --------
$ws.Columns("G").ColumnWidth = 16.830729166666668
$ws.Columns("H").ColumnWidth = 10.498697916666666
$ws.Columns("I").ColumnWidth = 11.998697916666666
$ws.Columns("J").ColumnWidth = 12.330729166666666
$ws.Columns("K").ColumnWidth = 20.830729166666668
$ws.Columns("L").ColumnWidth = 23.998697916666668

# --- Selection ---------------------------------------------------------------
$ws.Range("I12").Select() | Out-Null
